$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure D/E columns (price & volume text) are treated as text so values
# like "0.160" or "6.28" are not silently converted to numbers by Excel,
# then restore the default "Normal" style so no stray formatting is left
# behind (matching the original unstyled cells).
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "58.482.88"
$ws.Range("E2").Value = "  -4.47%  "

$ws.Range("D3").Value = "2.563.03"
$ws.Range("E3").Value = "  -4.05%  "

$ws.Range("E4").Value = "  +0.14%  "

$ws.Range("D5").Value = "508.18"
$ws.Range("E5").Value = "  -4.72%  "

$ws.Range("D6").Value = "145.83"
$ws.Range("E6").Value = "  -7.19%  "

$ws.Range("E7").Value = "  +0.20%  "

$ws.Range("D8").Value = "0.572"
$ws.Range("E8").Value = "  -3.38%  "

$ws.Range("D9").Value = "2.577.52"
$ws.Range("E9").Value = "  -4.17%  "

$ws.Range("D10").Value = "6.28"
$ws.Range("E10").Value = "  -5.19%  "

$ws.Range("E11").Value = "  -6.10%  "

$ws.Range("D12").Value = "0.335"
$ws.Range("E12").Value = "  -5.49%  "

$ws.Range("E13").Value = "  -0.92%  "

$ws.Range("D14").Value = "3.017.11"
$ws.Range("E14").Value = "  -3.56%  "

$ws.Range("D15").Value = "58.455.12"
$ws.Range("E15").Value = "  -4.48%  "

$ws.Range("D16").Value = "21.10"
$ws.Range("E16").Value = "  -4.93%  "

$ws.Range("E17").Value = "  -5.16%  "

$ws.Range("D18").Value = "2.580.16"
$ws.Range("E18").Value = "  -3.59%  "

$ws.Range("D19").Value = "4.54"
$ws.Range("E19").Value = "  -5.45%  "

$ws.Range("D20").Value = "344.10"
$ws.Range("E20").Value = "  -3.74%  "

$ws.Range("D21").Value = "10.27"
$ws.Range("E21").Value = "  -4.61%  "

$ws.Range("D22").Value = "6.06"
$ws.Range("E22").Value = "  -4.45%  "

$ws.Range("D23").Value = "0.998"
$ws.Range("E23").Value = "  -0.08%  "

$ws.Range("D24").Value = "60.70"
$ws.Range("E24").Value = "  -1.58%  "

$ws.Range("D25").Value = "0.417"
$ws.Range("E25").Value = "  -4.07%  "

$ws.Range("E26").Value = "  -0.21%  "

$ws.Range("B27").Value = "WrappedeETH"
$ws.Range("C27").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D27").Value = "2.688.09"
$ws.Range("E27").Value = "  -3.28%  "

$ws.Range("B28").Value = "Kaspa"
$ws.Range("C28").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D28").Value = "0.160"
$ws.Range("E28").Value = "  -5.42%  "

$ws.Range("D29").Value = "0.0₃0811"
$ws.Range("E29").Value = "  -6.73%  "

$ws.Range("D30").Value = "6.99"
$ws.Range("E30").Value = "  -6.24%  "

$ws.Range("E31").Value = "  +0.05%  "

$ws.Range("D32").Value = "6.07"
$ws.Range("E32").Value = "  -2.17%  "

$ws.Range("D33").Value = "18.74"
$ws.Range("E33").Value = "  -4.61%  "

$ws.Range("D34").Value = "149.81"
$ws.Range("E34").Value = "  -0.37%  "

$ws.Range("D35").Value = "1.54"
$ws.Range("E35").Value = "  -5.87%  "

$ws.Range("D36").Value = "0.945"
$ws.Range("E36").Value = "  +6.88%  "

$ws.Range("D37").Value = "3.97"
$ws.Range("E37").Value = "  -4.87%  "

$ws.Range("D38").Value = "1.13"
$ws.Range("E38").Value = "  -6.40%  "

$ws.Range("D39").Value = "0.854"
$ws.Range("E39").Value = "  -7.17%  "

$ws.Range("D40").Value = "36.03"
$ws.Range("E40").Value = "  -2.45%  "

$ws.Range("D41").Value = "292.16"
$ws.Range("E41").Value = "  -5.66%  "

$ws.Range("E42").Value = "  -7.39%  "

$ws.Range("D43").Value = "3.57"
$ws.Range("E43").Value = "  -6.87%  "

$ws.Range("D44").Value = "0.0993"
$ws.Range("E44").Value = "  -2.95%  "

$ws.Range("E45").Value = "  +0.06%  "

$ws.Range("D46").Value = "0.609"
$ws.Range("E46").Value = "  -6.66%  "

$ws.Range("D47").Value = "0.0536"
$ws.Range("E47").Value = "  -5.66%  "

$ws.Range("D48").Value = "19.07"
$ws.Range("E48").Value = "  -7.60%  "

$ws.Range("D49").Value = "10.26"
$ws.Range("E49").Value = "  -1.00%  "

$ws.Range("B50").Value = "RenderToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D50").Value = "4.69"
$ws.Range("E50").Value = "  -7.37%  "

$ws.Range("B51").Value = "VeChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D51").Value = "0.0227"
$ws.Range("E51").Value = "  -5.53%  "

# Restore default styling on the updated text cells (removes the temporary
# text-number-format override so the XML matches the original, unstyled cells).
$ws.Range("D2:E51").Style = "Normal"
